$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: Dept Loc "blore" -> Full address "NR peta"
$ws.Range("A5").Value = "blore"
$ws.Range("B5").Value = "NR peta"

# Update the selection to match the saved view state (B7)
$ws.Range("B7").Select()
